$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary - update aggregate metrics after closing Trade #12
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.99   # Current Capital
$wsSummary.Range("B4").Value = -0.01     # Total P&L $
$wsSummary.Range("B5").Value = -0.02     # Total P&L %
$wsSummary.Range("B6").Value = 12        # Total Trades
$wsSummary.Range("B7").Value = 5         # Winning Trades
$wsSummary.Range("B9").Value = 41.67     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: Strategy Status - update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.99   # Capital
$wsStatus.Range("D4").Value = 12      # Trades
$wsStatus.Range("E4").Value = -0.01   # P&L $
$wsStatus.Range("F4").Value = -0.01   # P&L %
$wsStatus.Range("G4").Value = 41.67   # Win Rate %

# ---------------------------------------------------------------------------
# Append new Trade #12 row (row 13) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$newTradeRow = @(12, "2026-02-17", "07:53:35", "MarketMaking", "DOWN", 0.21, 0.24, "CLOSED", 14.2857, 0.03, 99.99, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newTradeRow.Length; $col++) {
        $ws.Cells.Item(13, $col).Value = $newTradeRow[$col - 1]
    }
}
